# Scheduled data-refresh: update cached crafting-profit figures
# (currentAveragePrice / LevePrice* / LeveProfit* columns) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR market-data sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 2940
$ws.Range("J51").Value = 1940.909
$ws.Range("L51").Value = 1940.909
$ws.Range("N51").Value = -2908.909

# Row 103: Let Loose the Juice / Persimmon Tannin
$ws.Range("H103").Value = 1486.909
$ws.Range("I103").Value = 479.5
$ws.Range("K103").Value = 1438.5
$ws.Range("M103").Value = -852.5

# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 3704945
$ws.Range("I106").Value = 4167988.5
$ws.Range("K106").Value = 4167988.5
$ws.Range("M106").Value = -4167357.5

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 8001.6665
$ws.Range("I113").Value = 6005
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 6005
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -2751
$ws.Range("N113").Value = -15508

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2051.25
$ws.Range("I132").Value = 2238.1428
$ws.Range("K132").Value = 6714.428400000001
$ws.Range("M132").Value = -4184.428400000001

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 657.95
$ws.Range("J135").Value = 444.5
$ws.Range("L135").Value = 4000.5
$ws.Range("N135").Value = -9070.5

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3165.8572
$ws.Range("I137").Value = 3824.6667
$ws.Range("K137").Value = 11474.0001
$ws.Range("M137").Value = -8924.000100000001

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3220.5557
$ws.Range("I138").Value = 1510.8064
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 4532.4192
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 607.5807999999997
$ws.Range("N138").Value = -22280

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 26207452
$ws.Range("I32").Value = 29872266
$ws.Range("J32").Value = 8406931
$ws.Range("K32").Value = 29872266
$ws.Range("L32").Value = 8406931
$ws.Range("M32").Value = -29871979
$ws.Range("N32").Value = -8407505

# Row 37: Get Shirty / Steel Chainmail
$ws.Range("H37").Value = 37976.555
$ws.Range("J37").Value = 59994.2
$ws.Range("L37").Value = 59994.2
$ws.Range("N37").Value = -60540.2

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 6100
$ws.Range("I61").Value = 7050
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 7050
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -6838
$ws.Range("N61").Value = -4624

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 4684.091
$ws.Range("I122").Value = 4218.143
$ws.Range("J122").Value = 5499.5
$ws.Range("K122").Value = 12654.429
$ws.Range("L122").Value = 16498.5
$ws.Range("M122").Value = -10204.429
$ws.Range("N122").Value = -21398.5

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2181.625
$ws.Range("I132").Value = 2033.2894
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6099.8682
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3569.8682
$ws.Range("N132").Value = -20060

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 6100
$ws.Range("I136").Value = 7050
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 21150
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -18600
$ws.Range("N136").Value = -17700

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 2608.8
$ws.Range("I20").Value = 2799.8462
$ws.Range("J20").Value = 2254
$ws.Range("K20").Value = 2799.8462
$ws.Range("L20").Value = 2254
$ws.Range("M20").Value = -2552.8462
$ws.Range("N20").Value = -2748

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2644.0908
$ws.Range("I86").Value = 2399.1667
$ws.Range("J86").Value = 2938
$ws.Range("K86").Value = 2399.1667
$ws.Range("L86").Value = 2938
$ws.Range("M86").Value = -1276.1667
$ws.Range("N86").Value = -5184

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2644.0908
$ws.Range("I89").Value = 2399.1667
$ws.Range("J89").Value = 2938
$ws.Range("K89").Value = 11995.8335
$ws.Range("L89").Value = 14690
$ws.Range("M89").Value = -6379.833500000001
$ws.Range("N89").Value = -25922

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1631.5
$ws.Range("I31").Value = 950.0417
$ws.Range("K31").Value = 950.0417
$ws.Range("M31").Value = -655.0417

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1631.5
$ws.Range("I34").Value = 950.0417
$ws.Range("K34").Value = 950.0417
$ws.Range("M34").Value = -748.0417

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 3095.6296
$ws.Range("I58").Value = 1565.9333
$ws.Range("J58").Value = 5007.75
$ws.Range("K58").Value = 1565.9333
$ws.Range("L58").Value = 5007.75
$ws.Range("M58").Value = -1362.9333
$ws.Range("N58").Value = -5413.75

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 3441.8
$ws.Range("I99").Value = 3181
$ws.Range("J99").Value = 3833
$ws.Range("K99").Value = 3181
$ws.Range("L99").Value = 3833
$ws.Range("M99").Value = -1683
$ws.Range("N99").Value = -6829

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 3441.8
$ws.Range("I126").Value = 3181
$ws.Range("J126").Value = 3833
$ws.Range("K126").Value = 9543
$ws.Range("L126").Value = 11499
$ws.Range("M126").Value = -7073
$ws.Range("N126").Value = -16439

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 3095.6296
$ws.Range("I136").Value = 1565.9333
$ws.Range("J136").Value = 5007.75
$ws.Range("K136").Value = 4697.7999
$ws.Range("L136").Value = 15023.25
$ws.Range("M136").Value = -2147.7999
$ws.Range("N136").Value = -20123.25

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 182275310
$ws.Range("I4").Value = 155626990
$ws.Range("K4").Value = 466880970
$ws.Range("M4").Value = -466880858

# Row 33: Cooking with Gas / Chicken Stock
$ws.Range("H33").Value = 362.5
$ws.Range("I33").Value = 65.375
$ws.Range("J33").Value = 1551
$ws.Range("K33").Value = 392.25
$ws.Range("L33").Value = 9306
$ws.Range("M33").Value = -109.25
$ws.Range("N33").Value = -9872

# Row 44: No More Dumpster Diving / Knight's Bread
$ws.Range("H44").Value = 2500
$ws.Range("I44").Value = 2500
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 7500
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -7102
$ws.Range("N44").ClearContents()

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 1973.0435
$ws.Range("J68").Value = 2092.4666
$ws.Range("L68").Value = 6277.399800000001
$ws.Range("N68").Value = -7899.399800000001

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 1973.0435
$ws.Range("J71").Value = 2092.4666
$ws.Range("L71").Value = 18832.1994
$ws.Range("N71").Value = -26944.1994

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 1652.15
$ws.Range("I113").Value = 2111.111
$ws.Range("J113").Value = 1276.6364
$ws.Range("K113").Value = 6333.333
$ws.Range("L113").Value = 3829.9092
$ws.Range("M113").Value = -4163.333
$ws.Range("N113").Value = -8169.9092

# Row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 9000
$ws.Range("I133").Value = 2000
$ws.Range("K133").Value = 6000
$ws.Range("M133").Value = -940

$ws = $wb.Worksheets.Item("GSM")
# Row 26: Perk of Fiction / Coral Ring
$ws.Range("H26").Value = 40000
$ws.Range("J26").Value = 40000
$ws.Range("L26").Value = 40000
$ws.Range("N26").Value = -40560

# Row 50: Coral on My Mind / Red Coral Ring
$ws.Range("H50").Value = 40000
$ws.Range("J50").Value = 40000
$ws.Range("L50").Value = 40000
$ws.Range("N50").Value = -40996

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 461.91666
$ws.Range("I107").Value = 316.33334
$ws.Range("J107").Value = 607.5
$ws.Range("K107").Value = 316.33334
$ws.Range("L107").Value = 607.5
$ws.Range("M107").Value = 1603.66666
$ws.Range("N107").Value = -4447.5

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 19434.291
$ws.Range("I113").Value = 2359.5715
$ws.Range("J113").Value = 43338.9
$ws.Range("K113").Value = 2359.5715
$ws.Range("L113").Value = 43338.9
$ws.Range("M113").Value = -189.5715
$ws.Range("N113").Value = -47678.9

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3693.9443
$ws.Range("I126").Value = 3162.4
$ws.Range("K126").Value = 9487.200000000001
$ws.Range("M126").Value = -7017.200000000001

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 4750
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
# Row 9: From the Sands to the Stage / Leather Himantes
$ws.Range("H9").Value = 461.625
$ws.Range("I9").Value = 461.625
$ws.Range("K9").Value = 461.625
$ws.Range("M9").Value = -237.625

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 1499.8334
$ws.Range("I68").Value = 1399.8
$ws.Range("K68").Value = 1399.8
$ws.Range("M68").Value = -650.8

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 1499.8334
$ws.Range("I71").Value = 1399.8
$ws.Range("K71").Value = 6999
$ws.Range("M71").Value = -3255

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 35715412
$ws.Range("I93").Value = 52632508
$ws.Range("J93").Value = 1544.6666
$ws.Range("K93").Value = 52632508
$ws.Range("L93").Value = 1544.6666
$ws.Range("M93").Value = -52631260
$ws.Range("N93").Value = -4040.6666

# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 66629
$ws.Range("J133").Value = 66629
$ws.Range("L133").Value = 66629
$ws.Range("N133").Value = -71689

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 19953.264
$ws.Range("I136").Value = 22006.588
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 66019.764
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -63469.764
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
# Row 106: Cap It Off / Serge Knit Cap
$ws.Range("H106").Value = 89080
$ws.Range("I106").Value = 78342
$ws.Range("J106").Value = 92659.336
$ws.Range("K106").Value = 78342
$ws.Range("L106").Value = 92659.336
$ws.Range("M106").Value = -77080
$ws.Range("N106").Value = -95183.336

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 4510.1177
$ws.Range("I126").Value = 4178.1333
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 12534.3999
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -10064.3999
$ws.Range("N126").Value = -25940

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 35752.387
$ws.Range("I136").Value = 2998.25
$ws.Range("J136").Value = 95305.37
$ws.Range("K136").Value = 8994.75
$ws.Range("L136").Value = 285916.11
$ws.Range("M136").Value = -6444.75
$ws.Range("N136").Value = -291016.11

